$d = $word.ActiveDocument

# Move to the very end of the document body content
$endRange = $d.Content
$endRange.Collapse(0)

# Insert 3 blank paragraphs (matching lang en-GB formatting of surrounding paragraphs)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)

# Insert the final paragraph with the "Adding -> " text and the URL
$endRange.InsertAfter("Adding -> https://www.baeldung.com/spring-security-openid-connect")
